# "Generate Report for Handoff"
#
# The localization status workbook is regenerated once the files are
# handed off for translation: the Status / "Latest HO Xliff Generate
# Date" / "Latest Handoff Datetime" values move from "In Translation"
# to "Ready for handoff" (with refreshed timestamps), and the Status
# columns are widened to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet ------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-10-20 06:57:46"

# Widen the zh-cn / de-de status columns to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet -----------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-10-20 06:57:34"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet -----------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-10-20 06:57:46"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
